# Update fixed-cost results on row 2 of each year sheet with new values
# received from the server re-run of the "7_low_gas_demand" fix-cost case.

$wb = $excel.ActiveWorkbook

# Sheet "2025" (sheet1)
$ws = $wb.Worksheets.Item("2025")
$ws.Range("E2").Value = 29751.026632348
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 14142.32226446746
$ws.Range("L2").Value = 53029.15429284
$ws.Range("M2").Value = 10753.74970482998
$ws.Range("N2").Value = 7693.937482943782
$ws.Range("O2").Value = 7675.818959394448

# Sheet "2030" (sheet2)
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 4344.13718011466
$ws.Range("E2").Value = 57885.53306719843
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 28931.73319407337
$ws.Range("L2").Value = 90998.45501306068
$ws.Range("M2").Value = 20921.19698906325
$ws.Range("N2").Value = 10937.78653468269
$ws.Range("O2").Value = 9700.13305703022

# Sheets "2035", "2040", "2045", "2050" (sheet3-6) all receive identical
# updated values.
$years = @("2035", "2040", "2045", "2050")
foreach ($year in $years) {
    $ws = $wb.Worksheets.Item($year)
    $ws.Range("A2").Value = 1611.798000541006
    $ws.Range("B2").Value = 6181.041678721394
    $ws.Range("E2").Value = 71251.17376592311
    $ws.Range("G2").Value = 8095.925712661834
    $ws.Range("I2").Value = 45864.32315580232
    $ws.Range("L2").Value = 90998.45501306068
    $ws.Range("M2").Value = 25662.01124739233
    $ws.Range("N2").Value = 16047.23502266824
    $ws.Range("O2").Value = 15180.5186820528
}
